$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 0

# New data for rows 4-21 (Shot #, Label)
$data = @(
    @(3, 0),
    @(4, 1),
    @(5, 1),
    @(6, 1),
    @(7, 0),
    @(8, 0),
    @(9, 0),
    @(10, 0),
    @(11, 0),
    @(12, 0),
    @(13, 1),
    @(14, 0),
    @(15, 0),
    @(16, 1),
    @(17, 1),
    @(18, 0),
    @(19, 1),
    @(20, 0)
)

$row = 4
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Range("D21").Select()
